$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the per-row "time_taken" timestamps on the "data" sheet ---
$ws.Range("F2").Value = "2021-10-05 14:20:06.490053"
$ws.Range("F3").Value = "2021-10-05 14:20:06.490061"
$ws.Range("F4").Value = "2021-10-05 14:20:06.490065"
$ws.Range("F5").Value = "2021-10-05 14:20:06.490067"
$ws.Range("F6").Value = "2021-10-05 14:20:06.490070"
$ws.Range("F7").Value = "2021-10-05 14:20:06.490073"
$ws.Range("F8").Value = "2021-10-05 14:20:06.490075"
$ws.Range("F9").Value = "2021-10-05 14:20:06.490078"
$ws.Range("F10").Value = "2021-10-05 14:20:06.490081"
$ws.Range("F11").Value = "2021-10-05 14:20:06.490083"
$ws.Range("F12").Value = "2021-10-05 14:20:06.490086"
$ws.Range("F13").Value = "2021-10-05 14:20:06.490088"
$ws.Range("F14").Value = "2021-10-05 14:20:06.490091"
$ws.Range("F15").Value = "2021-10-05 14:20:06.490093"
$ws.Range("F16").Value = "2021-10-05 14:20:06.490096"
$ws.Range("F17").Value = "2021-10-05 14:20:06.490099"
$ws.Range("F18").Value = "2021-10-05 14:20:06.490101"
$ws.Range("F19").Value = "2021-10-05 14:20:06.490104"
$ws.Range("F20").Value = "2021-10-05 14:20:06.490107"
$ws.Range("F21").Value = "2021-10-05 14:20:06.490109"
$ws.Range("F22").Value = "2021-10-05 14:20:06.490112"
$ws.Range("F23").Value = "2021-10-05 14:20:06.490115"
$ws.Range("F24").Value = "2021-10-05 14:20:06.490117"
$ws.Range("F25").Value = "2021-10-05 14:20:06.490120"
$ws.Range("F26").Value = "2021-10-05 14:20:06.490123"
$ws.Range("F27").Value = "2021-10-05 14:20:06.490125"
$ws.Range("F28").Value = "2021-10-05 14:20:06.490128"
$ws.Range("F29").Value = "2021-10-05 14:20:06.490130"
$ws.Range("F30").Value = "2021-10-05 14:20:06.490133"
$ws.Range("F31").Value = "2021-10-05 14:20:06.490135"
$ws.Range("F32").Value = "2021-10-05 14:20:06.490138"
$ws.Range("F33").Value = "2021-10-05 14:20:06.490140"
$ws.Range("F34").Value = "2021-10-05 14:20:06.490143"
$ws.Range("F35").Value = "2021-10-05 14:20:06.490146"
$ws.Range("F36").Value = "2021-10-05 14:20:06.490148"
$ws.Range("F37").Value = "2021-10-05 14:20:06.490151"
$ws.Range("F38").Value = "2021-10-05 14:20:06.490153"
$ws.Range("F39").Value = "2021-10-05 14:20:06.490156"
$ws.Range("F40").Value = "2021-10-05 14:20:06.490158"
$ws.Range("F41").Value = "2021-10-05 14:20:06.490161"
$ws.Range("F42").Value = "2021-10-05 14:20:06.490164"
$ws.Range("F43").Value = "2021-10-05 14:20:06.490166"
$ws.Range("F44").Value = "2021-10-05 14:20:06.490169"
$ws.Range("F45").Value = "2021-10-05 14:20:06.490171"
$ws.Range("F46").Value = "2021-10-05 14:20:06.490174"
$ws.Range("F47").Value = "2021-10-05 14:20:06.490176"
$ws.Range("F48").Value = "2021-10-05 14:20:06.490179"
$ws.Range("F49").Value = "2021-10-05 14:20:06.490182"
$ws.Range("F50").Value = "2021-10-05 14:20:06.490184"
$ws.Range("F51").Value = "2021-10-05 14:20:06.490186"
$ws.Range("F52").Value = "2021-10-05 14:20:06.490189"
$ws.Range("F53").Value = "2021-10-05 14:20:06.490192"
$ws.Range("F54").Value = "2021-10-05 14:20:06.490194"
$ws.Range("F55").Value = "2021-10-05 14:20:06.490197"
$ws.Range("F56").Value = "2021-10-05 14:20:06.490199"
$ws.Range("F57").Value = "2021-10-05 14:20:06.490202"
$ws.Range("F58").Value = "2021-10-05 14:20:06.490205"
$ws.Range("F59").Value = "2021-10-05 14:20:06.490207"
$ws.Range("F60").Value = "2021-10-05 14:20:06.490210"
$ws.Range("F61").Value = "2021-10-05 14:20:06.490213"
$ws.Range("F62").Value = "2021-10-05 14:20:06.490215"
$ws.Range("F63").Value = "2021-10-05 14:20:06.490218"
$ws.Range("F64").Value = "2021-10-05 14:20:06.490220"
$ws.Range("F65").Value = "2021-10-05 14:20:06.490223"
$ws.Range("F66").Value = "2021-10-05 14:20:06.490226"
$ws.Range("F67").Value = "2021-10-05 14:20:06.490229"
$ws.Range("F68").Value = "2021-10-05 14:20:06.490232"
$ws.Range("F69").Value = "2021-10-05 14:20:06.490235"
$ws.Range("F70").Value = "2021-10-05 14:20:06.490237"
$ws.Range("F71").Value = "2021-10-05 14:20:06.490240"
$ws.Range("F72").Value = "2021-10-05 14:20:06.490242"
$ws.Range("F73").Value = "2021-10-05 14:20:06.490245"
$ws.Range("F74").Value = "2021-10-05 14:20:06.490248"
$ws.Range("F75").Value = "2021-10-05 14:20:06.490250"
$ws.Range("F76").Value = "2021-10-05 14:20:06.490253"
$ws.Range("F77").Value = "2021-10-05 14:20:06.490256"
$ws.Range("F78").Value = "2021-10-05 14:20:06.490260"
$ws.Range("F79").Value = "2021-10-05 14:20:06.490263"
$ws.Range("F80").Value = "2021-10-05 14:20:06.490266"
$ws.Range("F81").Value = "2021-10-05 14:20:06.490268"

# --- Add the new "metadata" tab, right after "data" ---
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Reuse the bold/bordered header style from the "data" sheet's header row
# (copy format only, so no new style entries are created).
$ws.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$meta.Range("G1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Ehlers Danlos syndromes"
$meta.Range("C2").Value = 53

# "2.63" looks numeric, so force text storage (like the source panel
# version string) then drop the direct formatting again so the cell keeps
# the default (unstyled) look.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "2.63"
$meta.Range("D2").ClearFormats()

$meta.Range("E2").Value = "2021-07-28T02:40:28.484710Z"
$meta.Range("F2").Value = "2021-10-05 14:20:06.486593"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/53/?format=json"

# Leave "data" as the active sheet/tab (matches activeTab="0" in the workbook)
$ws.Activate()
